$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new rows are treated as text (matches the existing sheet, where
# every data cell - including numeric- and date-looking values - is stored as text)
$ws.Range("A64:X80").NumberFormat = "@"

# Row 64
$ws.Range("A64").Value = "Fournisseur_05"
$ws.Range("B64").Value = "1020 – Banque 1"
$ws.Range("C64").Value = "34 jours"
$ws.Range("H64").Value = "9000"
$ws.Range("I64").Value = "9000"
$ws.Range("J64").Value = "9000"
$ws.Range("K64").Value = "9000"
$ws.Range("L64").Value = "9000"
$ws.Range("M64").Value = "1021 – Banque 2"
$ws.Range("N64").Value = "9000"
$ws.Range("O64").Value = "9000"
$ws.Range("P64").Value = "9000"
$ws.Range("Q64").Value = "81000000.00"
$ws.Range("R64").Value = "1100 – Débiteurs 1"
$ws.Range("S64").Value = "9000"
$ws.Range("T64").Value = "7.7"
$ws.Range("U64").Value = "9000"
$ws.Range("V64").Value = "81009000.00"
$ws.Range("X64").Value = "1015"

# Row 65
$ws.Range("A65").Value = "Fournisseur_03"
$ws.Range("B65").Value = "1010 – CCP 1"
$ws.Range("C65").Value = "32 jours"
$ws.Range("D65").Value = "2025-07-25"
$ws.Range("E65").Value = "2025-07-25"
$ws.Range("F65").Value = "2025-07-25"
$ws.Range("G65").Value = "2025-07"
$ws.Range("H65").Value = "8000"
$ws.Range("I65").Value = "8000"
$ws.Range("J65").Value = "8000"
$ws.Range("K65").Value = "8000"
$ws.Range("L65").Value = "8000"
$ws.Range("M65").Value = "1011 – CCP 2"
$ws.Range("N65").Value = "8000"
$ws.Range("O65").Value = "8000"
$ws.Range("P65").Value = "8000"
$ws.Range("Q65").Value = "64000000.00"
$ws.Range("R65").Value = "1020 – Banque 1"
$ws.Range("S65").Value = "8000"
$ws.Range("T65").Value = "7.7"
$ws.Range("U65").Value = "8000"
$ws.Range("V65").Value = "64008000.00"
$ws.Range("X65").Value = "1016"

# Row 66
$ws.Range("A66").Value = "Fournisseur_07"
$ws.Range("B66").Value = "1100 – Débiteurs 1"
$ws.Range("C66").Value = "36 jours"
$ws.Range("D66").Value = "2025-07-26"
$ws.Range("G66").Value = "2025-07"
$ws.Range("H66").Value = "20"
$ws.Range("I66").Value = "20"
$ws.Range("J66").Value = "20"
$ws.Range("K66").Value = "20"
$ws.Range("L66").Value = "20"
$ws.Range("M66").Value = "1101 – Débiteurs 2"
$ws.Range("N66").Value = "20"
$ws.Range("O66").Value = "20"
$ws.Range("P66").Value = "20"
$ws.Range("Q66").Value = "400.00"
$ws.Range("R66").Value = "1200 – Stock 1"
$ws.Range("S66").Value = "20"
$ws.Range("T66").Value = "7.7"
$ws.Range("U66").Value = "20"
$ws.Range("V66").Value = "420.00"
$ws.Range("X66").Value = "1017"

# Row 67
$ws.Range("A67").Value = "Fournisseur_04"
$ws.Range("B67").Value = "1011 – CCP 2"
$ws.Range("C67").Value = "33 jours"
$ws.Range("D67").Value = "2025-07-26"
$ws.Range("G67").Value = "2025-07"
$ws.Range("H67").Value = "3000"
$ws.Range("I67").Value = "3000"
$ws.Range("J67").Value = "3000"
$ws.Range("K67").Value = "3000"
$ws.Range("L67").Value = "3000"
$ws.Range("M67").Value = "1020 – Banque 1"
$ws.Range("N67").Value = "3000"
$ws.Range("O67").Value = "3000"
$ws.Range("P67").Value = "3000"
$ws.Range("Q67").Value = "9000000.00"
$ws.Range("R67").Value = "1021 – Banque 2"
$ws.Range("S67").Value = "3000"
$ws.Range("T67").Value = "7.7"
$ws.Range("U67").Value = "3000"
$ws.Range("V67").Value = "9003000.00"
$ws.Range("X67").Value = "1018"

# Row 68
$ws.Range("A68").Value = "Fournisseur_04"
$ws.Range("B68").Value = "1011 – CCP 2"
$ws.Range("C68").Value = "33 jours"
$ws.Range("H68").Value = "5000"
$ws.Range("J68").Value = "5000"
$ws.Range("M68").Value = "1020 – Banque 1"
$ws.Range("O68").Value = "5000"
$ws.Range("P68").Value = "1"
$ws.Range("Q68").Value = "5000.00"
$ws.Range("R68").Value = "1021 – Banque 2"
$ws.Range("T68").Value = "7.7"
$ws.Range("U68").Value = "5000"
$ws.Range("V68").Value = "10000.00"
$ws.Range("X68").Value = "1019"

# Row 69
$ws.Range("A69").Value = "Fournisseur_04"
$ws.Range("B69").Value = "1011 – CCP 2"
$ws.Range("C69").Value = "33 jours"
$ws.Range("D69").Value = "2025-07-26"
$ws.Range("H69").Value = "6000"
$ws.Range("J69").Value = "6000"
$ws.Range("M69").Value = "1020 – Banque 1"
$ws.Range("Q69").Value = "6000"
$ws.Range("R69").Value = "1021 – Banque 2"
$ws.Range("T69").Value = "7.7"
$ws.Range("U69").Value = "6000"
$ws.Range("V69").Value = "12000.00"
$ws.Range("X69").Value = "1020"

# Row 70
$ws.Range("A70").Value = "Fournisseur_10"
$ws.Range("B70").Value = "1201 – Stock2"
$ws.Range("C70").Value = "39 jours"
$ws.Range("G70").Value = "2025-07"
$ws.Range("H70").Value = "9000"
$ws.Range("M70").Value = "1300 – Actif transitoire 1"
$ws.Range("Q70").Value = "9000"
$ws.Range("R70").Value = "1301 – Actif transitoire 2"
$ws.Range("U70").Value = "9000"
$ws.Range("X70").Value = "1021"

# Row 71
$ws.Range("A71").Value = "Fournisseur_03"
$ws.Range("B71").Value = "1010 – CCP 1"
$ws.Range("C71").Value = "32 jours"
$ws.Range("H71").Value = "200"
$ws.Range("M71").Value = "1011 – CCP 2"
$ws.Range("Q71").Value = "200"
$ws.Range("R71").Value = "1020 – Banque 1"
$ws.Range("U71").Value = "200"
$ws.Range("X71").Value = "1022"

# Row 72
$ws.Range("A72").Value = "Fournisseur_06"
$ws.Range("B72").Value = "1021 – Banque 2"
$ws.Range("C72").Value = "35 jours"
$ws.Range("H72").Value = "300"
$ws.Range("M72").Value = "1100 – Débiteurs 1"
$ws.Range("Q72").Value = "300"
$ws.Range("R72").Value = "1101 – Débiteurs 2"
$ws.Range("U72").Value = "300"
$ws.Range("X72").Value = "1023"

# Row 73
$ws.Range("A73").Value = "Fournisseur_04"
$ws.Range("B73").Value = "1011 – CCP 2"
$ws.Range("C73").Value = "33 jours"
$ws.Range("H73").Value = "3"
$ws.Range("J73").Value = "3"
$ws.Range("M73").Value = "1020 – Banque 1"
$ws.Range("Q73").Value = "3"
$ws.Range("R73").Value = "1021 – Banque 2"
$ws.Range("U73").Value = "3"
$ws.Range("X73").Value = "1024"

# Row 74
$ws.Range("A74").Value = "Fournisseur_07"
$ws.Range("B74").Value = "1100 – Débiteurs 1"
$ws.Range("C74").Value = "36 jours"
$ws.Range("H74").Value = "15"
$ws.Range("M74").Value = "1101 – Débiteurs 2"
$ws.Range("Q74").Value = "10"
$ws.Range("R74").Value = "1200 – Stock 1"
$ws.Range("U74").Value = "5"
$ws.Range("X74").Value = "1025"

# Row 75
$ws.Range("A75").Value = "Fournisseur_10"
$ws.Range("B75").Value = "1201 – Stock2"
$ws.Range("C75").Value = "39 jours"
$ws.Range("G75").Value = "2025-07"
$ws.Range("M75").Value = "1300 – Actif transitoire 1"
$ws.Range("Q75").Value = "10"
$ws.Range("R75").Value = "1301 – Actif transitoire 2"
$ws.Range("U75").Value = "10"
$ws.Range("X75").Value = "1026"

# Row 76
$ws.Range("A76").Value = "Fournisseur_02"
$ws.Range("B76").Value = "1001 – Caisse 2"
$ws.Range("C76").Value = "31 jours"
$ws.Range("M76").Value = "1010 – CCP 1"
$ws.Range("Q76").Value = "2.00"
$ws.Range("R76").Value = "1011 – CCP 2"
$ws.Range("U76").Value = "2"
$ws.Range("X76").Value = "1027"

# Row 77
$ws.Range("A77").Value = "Fournisseur_05"
$ws.Range("B77").Value = "1020 – Banque 1"
$ws.Range("C77").Value = "34 jours"
$ws.Range("H77").Value = "4"
$ws.Range("M77").Value = "1021 – Banque 2"
$ws.Range("Q77").Value = "4"
$ws.Range("R77").Value = "1100 – Débiteurs 1"
$ws.Range("U77").Value = "4"
$ws.Range("X77").Value = "1028"

# Row 78
$ws.Range("B78").Value = "1011 – CCP 2"
$ws.Range("C78").Value = "33 jours"
$ws.Range("H78").Value = "5"
$ws.Range("M78").Value = "1020 – Banque 1"
$ws.Range("Q78").Value = "5"
$ws.Range("R78").Value = "1021 – Banque 2"
$ws.Range("U78").Value = "5"
$ws.Range("X78").Value = "1029"

# Row 79
$ws.Range("A79").Value = "Fournisseur_08"
$ws.Range("B79").Value = "1101 – Débiteurs 2"
$ws.Range("C79").Value = "37 jours"
$ws.Range("D79").Value = "2025-07-26"
$ws.Range("M79").Value = "1200 – Stock 1"
$ws.Range("Q79").Value = "8"
$ws.Range("R79").Value = "1201 – Stock2"
$ws.Range("U79").Value = "8"
$ws.Range("X79").Value = "1030"

# Row 80
$ws.Range("X80").Value = "1031"
